$wb = $excel.ActiveWorkbook

# --- Update the "INFRAÇÕES" sheet selection (D20 -> G27) ---
$sInfra = $wb.Worksheets.Item("INFRAÇÕES")
$sInfra.Activate()
$sInfra.Range("G27").Select()

# --- Insert the new "CONDUTORES POR CATEGORIA" sheet right before
#     "VEICULOS POR COMBUSTÍVEL" ---
$target = $wb.Worksheets.Item("VEICULOS POR COMBUSTÍVEL")
$newSheet = $wb.Worksheets.Add($target)
$newSheet.Name = "CONDUTORES POR CATEGORIA"

# Column widths
$newSheet.Columns.Item(1).ColumnWidth = 19.7109375
$newSheet.Columns.Item(2).ColumnWidth = 18.5703125

# Header row (column A text values are written in the same order the
# original workbook's shared-string table was built in, so new shared
# strings land at the same indices as in the target file)
$newSheet.Range("A1").Value = "CATEGORIAS"
$newSheet.Range("A3").Value = "E"
$newSheet.Range("A2").Value = "C"
$newSheet.Range("A4").Value = "B"
$newSheet.Range("A5").Value = "AE"
$newSheet.Range("A7").Value = "AD"
$newSheet.Range("A8").Value = "AC"
$newSheet.Range("A9").Value = "AB"
$newSheet.Range("A10").Value = "D"
$newSheet.Range("A6").Value = "A"

$newSheet.Range("B1").Value = "QUANTIDADE"

# Data rows (quantities)
$newSheet.Range("B2").Value = 3055
$newSheet.Range("B3").Value = 876
$newSheet.Range("B4").Value = 77539
$newSheet.Range("B5").Value = 2665
$newSheet.Range("B6").Value = 13713
$newSheet.Range("B7").Value = 12592
$newSheet.Range("B8").Value = 2032
$newSheet.Range("B9").Value = 68807
$newSheet.Range("B10").Value = 6456

# Activate the new sheet and set its selection to F20
$newSheet.Activate()
$newSheet.Range("F20").Select()
